$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.980.40'
$ws.Range('D3').Value = '2.670.21'
$ws.Range('E3').Value = '  -1.30%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.77'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '164.42'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.79%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D9').Value = '2.669.52'
$ws.Range('E9').Value = '  -1.16%  '
$ws.Range('E10').Value = '  +0.90%  '
$ws.Range('E11').Value = '  +0.68%  '
$ws.Range('E12').Value = '  -0.73%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.20'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.67'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.53%  '
$ws.Range('D15').Value = '3.156.17'
$ws.Range('E15').Value = '  -0.93%  '
$ws.Range('E16').Value = '  -2.61%  '
$ws.Range('D17').Value = '67.092.62'
$ws.Range('E17').Value = '  -2.03%  '
$ws.Range('D18').Value = '2.666.08'
$ws.Range('E18').Value = '  -2.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.63'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '360.32'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.87%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.49'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.55%  '
$ws.Range('E22').Value = '  -3.70%  '
$ws.Range('E23').Value = '  -2.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.03'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.98%  '
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '71.23'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.06'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.85%  '
$ws.Range('E29').Value = '  -2.20%  '
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '551.46'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.98'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.39'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.43%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.92'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.10%  '
$ws.Range('E35').Value = '  -2.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('E37').Value = '  -5.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.44'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.32%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '154.19'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.53%  '
$ws.Range('E40').Value = '  -2.42%  '
$ws.Range('E41').Value = '  -2.79%  '
$ws.Range('E42').Value = '  -4.63%  '
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('E44').Value = '  -4.85%  '
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.18'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.68%  '
$ws.Range('D47').Value = '0.0₆0297'
$ws.Range('E47').Value = '  -5.97%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.587'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.60%  '
$ws.Range('E49').Value = '  -3.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.83'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.73'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.76%  '
